# Update the "Förändrad" (last-changed) date in column C for every data row,
# and append the record's designation (column A) as the friendly display
# text of every HYPERLINK() formula found in the link columns (S..Y).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count - 1   # UsedRange starts at row 0 in this sheet
$firstLinkCol = 19   # column S
$lastLinkCol  = 25   # column Y

for ($r = 2; $r -le $lastRow; $r++) {

    $designation = $ws.Cells.Item($r, 1).Value2
    if ($designation -eq $null) { continue }

    # Column C holds the "changed" date serial - bump it from 45184 to 45186.
    $cCell = $ws.Cells.Item($r, 3)
    if ($cCell.Value2 -eq 45184) {
        $cCell.Value = 45186
    }

    # Walk the hyperlink columns and add the friendly name argument to any
    # HYPERLINK(...) formula that doesn't already have one.
    for ($col = $firstLinkCol; $col -le $lastLinkCol; $col++) {
        $cell = $ws.Cells.Item($r, $col)
        if ($cell.HasFormula) {
            $formula = $cell.Formula
            if ($formula -like 'HYPERLINK(*' -or $formula -like '=HYPERLINK(*') {
                if ($formula -notmatch ',') {
                    $trimmed = $formula.Substring(0, $formula.Length - 1)
                    $cell.Formula = $trimmed + ', "' + $designation + '")'
                }
            }
        }
    }
}
